$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Laatst bijgewerkt" timestamp cell (A2)
$ws.Range("A2").Value = "Laatst bijgewerkt: 2025-09-06 17:48:36"

# 2. Remove existing hyperlinks so M/N columns can be safely rewritten
$ws.Hyperlinks.Delete()

# 3. Rewrite the full data table (rows 3-16)
$rowData = @{
  3 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 9.5"; "E"="toto"; "F"=2.55; "G"="minder dan 9.5"; "H"="jacks"; "I"=1.77; "J"="1=61, 2=89"; "K"="€5.55"; "L"=4.29; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://jacks.nl/sports/event/1023224918#event/1023224918" }
  4 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 8.5"; "E"="toto"; "F"=1.9; "G"="minder dan 8.5"; "H"="jacks"; "I"=2.3; "J"="1=82, 2=68"; "K"="€5.8"; "L"=3.89; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://jacks.nl/sports/event/1023224918#event/1023224918" }
  5 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="armenië"; "D"="meer dan 2.5"; "E"="toto"; "F"=2.65; "G"="minder dan 2.5"; "H"="jacks"; "I"=1.71; "J"="1=59, 2=91"; "K"="€5.61"; "L"=3.78; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://jacks.nl/sports/event/1023224918#event/1023224918" }
  6 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 8.5"; "E"="toto"; "F"=1.9; "G"="minder dan 8.5"; "H"="kambi"; "I"=2.28; "J"="1=82, 2=68"; "K"="€5.04"; "L"=3.51; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://jacks.nl/sports/event/1023224918#event/1023224918" }
  7 = @{ "A"="Duitsland vs Noord-Ierland"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 10.5"; "E"="toto"; "F"=2.5; "G"="minder dan 10.5"; "H"="starcasino"; "I"=1.76; "J"="1=62, 2=88"; "K"="€4.88"; "L"=3.18; "M"="https://sport.toto.nl/wedden/wedstrijd/8778584"; "N"="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394" }
  8 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 7.5"; "E"="toto"; "F"=3.45; "G"="minder dan 10.5"; "H"="jacks"; "I"=1.46; "J"="1=45, 2=105"; "K"="€3.3"; "L"=2.52; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://jacks.nl/sports/event/1023224918#event/1023224918" }
  9 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="armenië"; "D"="meer dan 1.5"; "E"="toto"; "F"=1.55; "G"="minder dan 1.5"; "H"="jacks"; "I"=3; "J"="1=99, 2=51"; "K"="€3.0"; "L"=2.15; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://jacks.nl/sports/event/1023224918#event/1023224918" }
  10 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 7.5"; "E"="toto"; "F"=1.5; "G"="minder dan 7.5"; "H"="jacks"; "I"=3.2; "J"="1=102, 2=48"; "K"="€3.0"; "L"=2.08; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://jacks.nl/sports/event/1023224918#event/1023224918" }
  11 = @{ "A"="Litouwen vs Nederland"; "B"="totaal aantal schoten op doel"; "C"="nederland"; "D"="meer dan 8.5"; "E"="toto"; "F"=2.45; "G"="minder dan 8.5"; "H"="jacks"; "I"=1.74; "J"="1=62, 2=88"; "K"="€1.9"; "L"=1.71; "M"="https://sport.toto.nl/wedden/wedstrijd/8706282"; "N"="https://jacks.nl/sports/event/1023224945#event/1023224945" }
  12 = @{ "A"="Armenië vs Portugal"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 7.5"; "E"="toto"; "F"=1.5; "G"="minder dan 7.5"; "H"="kambi"; "I"=3.15; "J"="1=102, 2=48"; "K"="€1.2"; "L"=1.59; "M"="https://sport.toto.nl/wedden/wedstrijd/8590793"; "N"="https://www.unibet.nl/betting/sports/event/1023224918?coupon=single%7C3856344734%7C3.15%7Creplace" }
  13 = @{ "A"="Engeland vs Andorra"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 11.5"; "E"="toto"; "F"=2.6; "G"="minder dan 11.5"; "H"="starcasino"; "I"=1.6667; "J"="1=59, 2=91"; "K"="€1.67"; "L"=1.54; "M"="https://sport.toto.nl/wedden/wedstrijd/8668833"; "N"="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=13549445" }
  14 = @{ "A"="Duitsland vs Noord-Ierland"; "B"="totaal aantal schoten op doel"; "C"="wedstrijd"; "D"="meer dan 10.5"; "E"="toto"; "F"=2.5; "G"="minder dan 10.5"; "H"="onecasino"; "I"=1.71; "J"="1=61, 2=89"; "K"="€2.19"; "L"=1.52; "M"="https://sport.toto.nl/wedden/wedstrijd/8778584"; "N"="https://sports.onecasino.nl/#/event/10028349" }
  15 = @{ "A"="Litouwen vs Nederland"; "B"="totaal aantal schoten op doel"; "C"="nederland"; "D"="meer dan 7.5"; "E"="toto"; "F"=1.85; "G"="minder dan 7.5"; "H"="jacks"; "I"=2.25; "J"="1=82, 2=68"; "K"="€1.7"; "L"=1.5; "M"="https://sport.toto.nl/wedden/wedstrijd/8706282"; "N"="https://jacks.nl/sports/event/1023224945#event/1023224945" }
  16 = @{ "A"="Litouwen vs Nederland"; "B"="totaal aantal schoten op doel"; "C"="nederland"; "D"="meer dan 9.5"; "E"="toto"; "F"=3.5; "G"="minder dan 9.5"; "H"="jacks"; "I"=1.43; "J"="1=44, 2=106"; "K"="€1.58"; "L"=1.5; "M"="https://sport.toto.nl/wedden/wedstrijd/8706282"; "N"="https://jacks.nl/sports/event/1023224945#event/1023224945" }
}

$colIndex = @{ "A"=1; "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "G"=7; "H"=8; "I"=9; "J"=10; "K"=11; "L"=12; "M"=13; "N"=14 }

foreach ($r in $rowData.Keys) {
    $data = $rowData[$r]
    foreach ($col in $data.Keys) {
        $c = $colIndex[$col]
        $ws.Cells.Item($r, $c).Value = $data[$col]
    }
}

# 4. Ensure M/N columns (Link 1 / Link 2) use the Hyperlink style for the new rows
$ws.Range("M14:N16").Style = "Hyperlink"

# 5. Re-create hyperlinks for M3:N16 from the (address, sub-address) pairs below
$linkData = @{
  3 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224918"; sub="event/1023224918"} }
  4 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224918"; sub="event/1023224918"} }
  5 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224918"; sub="event/1023224918"} }
  6 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224918"; sub="event/1023224918"} }
  7 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8778584"; sub=None}; N=@{addr="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394"; sub=None} }
  8 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224918"; sub="event/1023224918"} }
  9 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224918"; sub="event/1023224918"} }
  10 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224918"; sub="event/1023224918"} }
  11 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8706282"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224945"; sub="event/1023224945"} }
  12 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8590793"; sub=None}; N=@{addr="https://www.unibet.nl/betting/sports/event/1023224918?coupon=single%7C3856344734%7C3.15%7Creplace"; sub=None} }
  13 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8668833"; sub=None}; N=@{addr="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=13549445"; sub=None} }
  14 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8778584"; sub=None}; N=@{addr="https://sports.onecasino.nl/"; sub="/event/10028349"} }
  15 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8706282"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224945"; sub="event/1023224945"} }
  16 = @{ M=@{addr="https://sport.toto.nl/wedden/wedstrijd/8706282"; sub=None}; N=@{addr="https://jacks.nl/sports/event/1023224945"; sub="event/1023224945"} }
}

foreach ($r in $linkData.Keys) {
    $entry = $linkData[$r]
    foreach ($col in @("M","N")) {
        $info = $entry[$col]
        $cell = $ws.Cells.Item($r, $colIndex[$col])
        if ($info.sub -eq $null) {
            $ws.Hyperlinks.Add($cell, $info.addr)
        } else {
            $ws.Hyperlinks.Add($cell, $info.addr, $info.sub)
        }
    }
}

Write-Host "Advies tabblad updated"